# Translate validate_required_fields_use_case.dart strings to English and
# add them to the "Domain" sheet (sheet2), right after the
# validate_item_name_use_case.dart section.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Clear out the old empty buffer rows (20:23) completely so we start
#    from a clean, unstyled slate below row 19 (same as the genuinely
#    blank separator rows 11 and 16 elsewhere in the sheet).
# ---------------------------------------------------------------------
$ws.Range("A20:B23").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2. Row 20: new English/Vietnamese pair for the "already taken" message
#    (reuses the existing English string from row 18, paired with a new
#    Vietnamese translation). Keep the tall wrapped-text row height.
# ---------------------------------------------------------------------
$ws.Range("A18:B18").Copy()
$ws.Range("A20:B20").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A20").Value = $ws.Range("A18").Value2
$ws.Range("B20").Value = "… của món đồ ... đã được sử dụng. Bạn vui lòng nhập tên khác. Có thể thêm số vào sau tên đồ vật (ví dụ: Áo 1, Áo 2,... để phân biệt)."
$ws.Rows.Item(20).RowHeight = 47.25

# Row 21 is intentionally left untouched -> stays a blank separator row
# (no row 21 is written out), matching rows 11 and 16.

# ---------------------------------------------------------------------
# 3. Row 22: new section header "validate_required_fields_use_case.dart"
#    Merge first, then copy the header formatting from row 17, then set
#    the text.
# ---------------------------------------------------------------------
$ws.Range("A22:B22").Merge()

$ws.Range("A17:B17").Copy()
$ws.Range("A22:B22").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A22").Value = "validate_required_fields_use_case.dart"
$ws.Rows.Item(22).RowHeight = 15.75

# ---------------------------------------------------------------------
# 4. Rows 23-29: copy the plain content-row formatting down, then fill
#    in the English/Vietnamese translation pairs. Rows 28-29 stay empty
#    buffer rows (matching the pattern used elsewhere in the sheet).
# ---------------------------------------------------------------------
$ws.Range("A18:B18").Copy()
$ws.Range("A23:B29").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

foreach ($r in 23..29) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

$ws.Range("A23").Value = "Please enter item name"
$ws.Range("B23").Value = "Vui lòng nhập tên đồ vật"

$ws.Range("A24").Value = "Please select a closet"
$ws.Range("B24").Value = "Vui lòng chọn tủ đồ"

$ws.Range("A25").Value = "Please select a category"
$ws.Range("B25").Value = "Vui lòng chọn danh mục"

$ws.Range("A26").Value = "Please enter a name for Item"
$ws.Range("B26").Value = "Vui lòng nhập tên cho đồ vật số"

$ws.Range("A27").Value = "Please select a closet for Item"
$ws.Range("B27").Value = "Vui lòng chọn tủ quần áo cho đồ vật số"

# Rows 28 and 29 remain empty (blank trailing buffer rows), as in the
# original template.

# ---------------------------------------------------------------------
# 5. Update the view selection to match the edited sheet (last edited
#    cell A28), and keep the sheet active.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A28").Select()
